$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. "Please find attached the manuscript" -> "W are pleased to submit our manuscript"
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "Please find attached the manuscript", $true, $false, $false, $false, $false,
    $true, 1, $false, "W are pleased to submit our manuscript", 2)

# ---------------------------------------------------------------------
# 2. "which we are submitting for consideration as a Letter in" -> "for consideration as a Letter in"
# ---------------------------------------------------------------------
$d.Content.Find.Execute(
    "which we are submitting for consideration as a Letter in", $true, $false, $false, $false, $false,
    $true, 1, $false, "for consideration as a Letter in", 2)

# ---------------------------------------------------------------------
# 3. Rewrite the "When a species goes extinct..." paragraph.
# ---------------------------------------------------------------------
$old3 = "When a species goes extinct, more than a species is lost. That species interacted with other species in its community, likely in ways that were important or even crucial to others" + [char]8217 + " survival. As the world continues to lose species at an alarming rate, it has become increasingly imperative to aid the recovery of lost interactions and component biodiversity through ecological restoration. We know little, however, about how to re-assemble interacting communities through restoration, or the process of ecological network assembly more generally."
$new3 = "Our research deals with two fundamental aspects of ecological theory: unerstanding how species-rich communities assemble, and how these assemblages change through time. As the world continues to lose species at an alarming rate, it has become increasingly imperative to aid the recovery of lost interactions and component biodiversity through ecological restoration. When a species goes extinct, not only a species is lost, but also its interactions. We know little, however, about how to re-assemble interacting communities through restoration, or the process of ecological network assembly more generally."
$d.Content.Find.Execute($old3, $true, $false, $false, $false, $false, $true, 1, $false, $new3, 2)

# ---------------------------------------------------------------------
# 4. Rewrite the "pollinator visitation records..." run, then append a
#    new " Ecology Letters." tail (with "Ecology Letters" italicised)
#    after it, inside the same paragraph.
# ---------------------------------------------------------------------
$old4 = "pollinator visitation records, we explore the assembly of plant-pollinator communities at native plant restorations in the Central Valley of California. Employing newly developed methods for examining temporal changes in networks, we find that species are highly dynamic in their network position, causing community assembly to be punctuated by major interaction reorganizations. The most persistent and generalized species were also the most variable in their network positions, contrary to what is expected through preferential attachment theory " + [char]8212 + " an assembly theory otherwise well-supported in the network literature. Our study is the first long-term study on the temporal assembly of ecological networks. It also challenges the hypothesis that mutualistic systems assemble through preferential attachment (Bascompte and Stouffer, 2009)."
$new4 = "pollinator visitation records, we explore the assembly of plant-pollinator communities at native plant restorations in the Central Valley of California. For the first time in the ecological litterature we employ a newly developed method to examine the temporal changes in networks. Among other things, we find that species are highly dynamic in their network position, causing community assembly to be punctuated by major interaction reorganizations. The most persistent and generalized species are also the most variable in their network positions, contrary to what is expected through preferential attachment theory " + [char]8212 + " an assembly theory otherwise well-supported in the network literature. Our study is the first long-term study on the temporal assembly of ecological networks. Our results are compelling and provide empirical evidence that widen our understading on how communities assembly and how species interactions changes through time. Furthermore, our results also contribute to the knowledge of how communities will be able to maintain function in the face of species extinction. And finally, our results challenge the view that communities assemble through preferential attachment. We believe that these exciting results that link three major ecological fields (interaction networks, community dinamics and restauration ecology) that will be of broad interest to the readership of"

$fr = $d.Content
$fr.Find.ClearFormatting()
$fr.Find.Text = $old4
$fr.Find.Execute()
$fr.Text = $new4
$pos = $fr.End

# " "
$d.Range($pos, $pos).InsertAfter(" ")
$pos = $pos + 1

# "Ecology Letters" (italic)
$ecoText = "Ecology Letters"
$d.Range($pos, $pos).InsertAfter($ecoText)
$ecoStart = $pos
$pos = $pos + $ecoText.Length
$ecoRange = $d.Range($ecoStart, $pos)
$ecoRange.Italic = 1

# "."
$d.Range($pos, $pos).InsertAfter(".")
$pos = $pos + 1

# ---------------------------------------------------------------------
# 5. Shorten the final "Our manuscript is original..." paragraph down to
#    just the closing "Thank you..." sentence.
# ---------------------------------------------------------------------
$old5 = "Our manuscript is original and was carried out fully by the authors. All authors agree with the contents of the manuscript. This manuscript is not published, nor is it in consideration for publication elsewhere. All research not of the authors" + [char]8217 + " is fully acknowledged. The authors declare no conflict of interest. All appropriate ethical standards were followed. Thank you for reviewing our manuscript and we hope you will find it suitable for publication."
$new5 = "Thank you for reviewing our manuscript and we hope you will find it suitable for publication."
$d.Content.Find.Execute($old5, $true, $false, $false, $false, $false, $true, 1, $false, $new5, 2)

Write-Output $d.Content.Text
